# Update tests. Model was not picking up the correct soil depth value
# from test configs.
#
# Row 5  (B5 = "Texture")     : C5:Q5 changes from "Silt loam" -> "Silt"
# Row 7  (B7 = "SampleDepth") : C7:Q7 changes from "0-30cm"    -> "Top30cm"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C5:Q5").Value = "Silt"
$ws.Range("C7:Q7").Value = "Top30cm"
